# Refresh the cryptos list with the latest values from the GitHub Actions run.
#
# Most rows only get updated Price (column D) and Volume(1h) (column E)
# figures. Two pairs of rows also swapped ranking positions between runs
# (EthereumClassic/Fetch.AI at rows 34-35, and Hedera/Mantle at rows 47-48),
# so their Coin (B) and Link (C) columns are rewritten too.
#
# Column D holds prices stored as literal text (e.g. "64.062.11",
# "4.91") in the source file, not numbers - a leading apostrophe is used
# below to force Excel to keep the new values as text as well, instead of
# re-parsing them as numbers (which would drop trailing zeros like
# "4.90" -> 4.9 or mis-handle the "." thousands separators).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.440.82"
$ws.Range("E2").Value = "  +6.79%  "
$ws.Range("D3").Value = "'2.739.73"
$ws.Range("E3").Value = "  +4.81%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'593.33"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").Value = "'152.48"
$ws.Range("E6").Value = "  +6.47%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  +2.42%  "
$ws.Range("D9").Value = "'2.779.10"
$ws.Range("E9").Value = "  +5.96%  "
$ws.Range("E10").Value = "  +3.85%  "
$ws.Range("E11").Value = "  +7.74%  "
$ws.Range("E12").Value = "  +3.69%  "
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").Value = "'3.236.58"
$ws.Range("E14").Value = "  +5.31%  "
$ws.Range("E15").Value = "  +7.43%  "
$ws.Range("D16").Value = "'64.230.16"
$ws.Range("E16").Value = "  +6.47%  "
$ws.Range("E17").Value = "  +9.11%  "
$ws.Range("D18").Value = "'2.769.32"
$ws.Range("E18").Value = "  +5.82%  "
$ws.Range("D19").Value = "'12.10"
$ws.Range("E19").Value = "  +5.59%  "
$ws.Range("D20").Value = "'4.90"
$ws.Range("E20").Value = "  +5.60%  "
$ws.Range("D21").Value = "'365.92"
$ws.Range("E21").Value = "  +5.61%  "
$ws.Range("D22").Value = "'7.03"
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").Value = "'0.994"
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("D25").Value = "'66.38"
$ws.Range("E25").Value = "  +4.16%  "
$ws.Range("D26").Value = "'0.168"
$ws.Range("E26").Value = "  +4.78%  "
$ws.Range("D27").Value = "'8.65"
$ws.Range("E27").Value = "  +7.61%  "
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").Value = "'0.0₃0910"
$ws.Range("E29").Value = "  +14.04%  "
$ws.Range("E30").Value = "  +4.99%  "
$ws.Range("E31").Value = "  +10.29%  "
$ws.Range("D32").Value = "'172.05"
$ws.Range("E32").Value = "  +2.03%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.19"
$ws.Range("E34").Value = "  +16.00%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'20.69"
$ws.Range("E35").Value = "  +6.19%  "
$ws.Range("D36").Value = "'4.80"
$ws.Range("E36").Value = "  +11.62%  "
$ws.Range("E37").Value = "  +9.65%  "
$ws.Range("E38").Value = "  +9.53%  "
$ws.Range("D39").Value = "'1.01"
$ws.Range("E39").Value = "  +19.26%  "
$ws.Range("D40").Value = "'349.20"
$ws.Range("E40").Value = "  +9.50%  "
$ws.Range("E41").Value = "  +7.97%  "
$ws.Range("D42").Value = "'39.44"
$ws.Range("E42").Value = "  +2.70%  "
$ws.Range("D43").Value = "'5.66"
$ws.Range("E43").Value = "  +12.71%  "
$ws.Range("D44").Value = "'22.28"
$ws.Range("E44").Value = "  +11.68%  "
$ws.Range("D45").Value = "'143.93"
$ws.Range("E45").Value = "  +6.13%  "
$ws.Range("D46").Value = "'22.12"
$ws.Range("E46").Value = "  +10.33%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.651"
$ws.Range("E47").Value = "  +7.15%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "'0.0593"
$ws.Range("E48").Value = "  +7.65%  "
$ws.Range("E49").Value = "  +7.49%  "
$ws.Range("E50").Value = "  +2.57%  "
$ws.Range("D51").Value = "'2.168.29"
$ws.Range("E51").Value = "  +7.11%  "
